$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 9.108069666666667
$ws.Cells.Item(2, 8).Value = 27.324209
$ws.Cells.Item(2, 9).Value = 0.00155006418458712
$ws.Cells.Item(2, 10).Value = 0.00155006418458712
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 417.3096593333333
$ws.Cells.Item(2, 14).Value = 1251.928978
$ws.Cells.Item(2, 15).Value = 0.8277575129472603
$ws.Cells.Item(2, 16).Value = 0.8277575129472603
$ws.Cells.Item(2, 17).Value = 3800.885449780933
$ws.Cells.Item(2, 18).Value = 34207.9690480284
$ws.Cells.Item(2, 19).Value = 0.001283077274342457
$ws.Cells.Item(2, 20).Value = 0.001283077274342457
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 9.108069666666667
$ws.Cells.Item(3, 8).Value = 27.324209
$ws.Cells.Item(3, 9).Value = 0.00155006418458712
$ws.Cells.Item(3, 10).Value = 0.00155006418458712
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.194815333333333
$ws.Cells.Item(3, 14).Value = 3.584446
$ws.Cells.Item(3, 15).Value = 0.002369984366839822
$ws.Cells.Item(3, 16).Value = 0.002369984366839822
$ws.Cells.Item(3, 17).Value = 10.88246129480156
$ws.Cells.Item(3, 18).Value = 97.942151653214
$ws.Cells.Item(3, 19).Value = 0.00000367362788506979
$ws.Cells.Item(3, 20).Value = 0.000003673627885069789
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 9.108069666666667
$ws.Cells.Item(4, 8).Value = 27.324209
$ws.Cells.Item(4, 9).Value = 0.00155006418458712
$ws.Cells.Item(4, 10).Value = 0.00155006418458712
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 58.96764933333333
$ws.Cells.Item(4, 14).Value = 176.902948
$ws.Cells.Item(4, 15).Value = 0.1169656960121252
$ws.Cells.Item(4, 16).Value = 0.1169656960121252
$ws.Cells.Item(4, 17).Value = 537.0814582075702
$ws.Cells.Item(4, 18).Value = 4833.733123868132
$ws.Cells.Item(4, 19).Value = 0.0001813043362136997
$ws.Cells.Item(4, 20).Value = 0.0001813043362136997
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.108069666666667
$ws.Cells.Item(5, 8).Value = 27.324209
$ws.Cells.Item(5, 9).Value = 0.00155006418458712
$ws.Cells.Item(5, 10).Value = 0.00155006418458712
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 26.67269233333333
$ws.Cells.Item(5, 14).Value = 80.018077
$ws.Cells.Item(5, 15).Value = 0.05290680667377473
$ws.Cells.Item(5, 16).Value = 0.05290680667377473
$ws.Cells.Item(5, 17).Value = 242.9367399695659
$ws.Cells.Item(5, 18).Value = 2186.430659726093
$ws.Cells.Item(5, 19).Value = 0.000082008946145893
$ws.Cells.Item(5, 20).Value = 0.000082008946145893
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5771.873535333333
$ws.Cells.Item(6, 8).Value = 17315.620606
$ws.Cells.Item(6, 9).Value = 0.9822909543423312
$ws.Cells.Item(6, 10).Value = 0.9822909543423313
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 417.3096593333333
$ws.Cells.Item(6, 14).Value = 1251.928978
$ws.Cells.Item(6, 15).Value = 0.8277575129472603
$ws.Cells.Item(6, 16).Value = 0.8277575129472603
$ws.Cells.Item(6, 17).Value = 2408658.578745035
$ws.Cells.Item(6, 18).Value = 21677927.20870532
$ws.Cells.Item(6, 19).Value = 0.8130987173569989
$ws.Cells.Item(6, 20).Value = 0.8130987173569989
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5771.873535333333
$ws.Cells.Item(7, 8).Value = 17315.620606
$ws.Cells.Item(7, 9).Value = 0.9822909543423312
$ws.Cells.Item(7, 10).Value = 0.9822909543423313
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.194815333333333
$ws.Cells.Item(7, 14).Value = 3.584446
$ws.Cells.Item(7, 15).Value = 0.002369984366839822
$ws.Cells.Item(7, 16).Value = 0.002369984366839822
$ws.Cells.Item(7, 17).Value = 6896.323002077142
$ws.Cells.Item(7, 18).Value = 62066.90701869428
$ws.Cells.Item(7, 19).Value = 0.002328014205479494
$ws.Cells.Item(7, 20).Value = 0.002328014205479494
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5771.873535333333
$ws.Cells.Item(8, 8).Value = 17315.620606
$ws.Cells.Item(8, 9).Value = 0.9822909543423312
$ws.Cells.Item(8, 10).Value = 0.9822909543423313
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 58.96764933333333
$ws.Cells.Item(8, 14).Value = 176.902948
$ws.Cells.Item(8, 15).Value = 0.1169656960121252
$ws.Cells.Item(8, 16).Value = 0.1169656960121252
$ws.Cells.Item(8, 17).Value = 340353.814627883
$ws.Cells.Item(8, 18).Value = 3063184.331650947
$ws.Cells.Item(8, 19).Value = 0.1148943451610654
$ws.Cells.Item(8, 20).Value = 0.1148943451610654
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5771.873535333333
$ws.Cells.Item(9, 8).Value = 17315.620606
$ws.Cells.Item(9, 9).Value = 0.9822909543423312
$ws.Cells.Item(9, 10).Value = 0.9822909543423313
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 26.67269233333333
$ws.Cells.Item(9, 14).Value = 80.018077
$ws.Cells.Item(9, 15).Value = 0.05290680667377473
$ws.Cells.Item(9, 16).Value = 0.05290680667377473
$ws.Cells.Item(9, 17).Value = 153951.406994855
$ws.Cells.Item(9, 18).Value = 1385562.662953695
$ws.Cells.Item(9, 19).Value = 0.05196987761878739
$ws.Cells.Item(9, 20).Value = 0.0519698776187874
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.272029666666667
$ws.Cells.Item(10, 8).Value = 3.816089
$ws.Cells.Item(10, 9).Value = 0.0002164813950916887
$ws.Cells.Item(10, 10).Value = 0.0002164813950916887
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 417.3096593333333
$ws.Cells.Item(10, 14).Value = 1251.928978
$ws.Cells.Item(10, 15).Value = 0.8277575129472603
$ws.Cells.Item(10, 16).Value = 0.8277575129472603
$ws.Cells.Item(10, 17).Value = 530.8302668585601
$ws.Cells.Item(10, 18).Value = 4777.472401727042
$ws.Cells.Item(10, 19).Value = 0.0001791941012004495
$ws.Cells.Item(10, 20).Value = 0.0001791941012004495
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.272029666666667
$ws.Cells.Item(11, 8).Value = 3.816089
$ws.Cells.Item(11, 9).Value = 0.0002164813950916887
$ws.Cells.Item(11, 10).Value = 0.0002164813950916887
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.194815333333333
$ws.Cells.Item(11, 14).Value = 3.584446
$ws.Cells.Item(11, 15).Value = 0.002369984366839822
$ws.Cells.Item(11, 16).Value = 0.002369984366839822
$ws.Cells.Item(11, 17).Value = 1.519840550188222
$ws.Cells.Item(11, 18).Value = 13.678564951694
$ws.Cells.Item(11, 19).Value = 0.0000005130575220789772
$ws.Cells.Item(11, 20).Value = 0.0000005130575220789772
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.272029666666667
$ws.Cells.Item(12, 8).Value = 3.816089
$ws.Cells.Item(12, 9).Value = 0.0002164813950916887
$ws.Cells.Item(12, 10).Value = 0.0002164813950916887
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 58.96764933333333
$ws.Cells.Item(12, 14).Value = 176.902948
$ws.Cells.Item(12, 15).Value = 0.1169656960121252
$ws.Cells.Item(12, 16).Value = 0.1169656960121252
$ws.Cells.Item(12, 17).Value = 75.00859932559689
$ws.Cells.Item(12, 18).Value = 675.077393930372
$ws.Cells.Item(12, 19).Value = 0.00002532089705057522
$ws.Cells.Item(12, 20).Value = 0.00002532089705057523
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.272029666666667
$ws.Cells.Item(13, 8).Value = 3.816089
$ws.Cells.Item(13, 9).Value = 0.0002164813950916887
$ws.Cells.Item(13, 10).Value = 0.0002164813950916887
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 26.67269233333333
$ws.Cells.Item(13, 14).Value = 80.018077
$ws.Cells.Item(13, 15).Value = 0.05290680667377473
$ws.Cells.Item(13, 16).Value = 0.05290680667377473
$ws.Cells.Item(13, 17).Value = 33.92845593787256
$ws.Cells.Item(13, 18).Value = 305.356103440853
$ws.Cells.Item(13, 19).Value = 0.00001145333931858502
$ws.Cells.Item(13, 20).Value = 0.00001145333931858502
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 93.67702500000001
$ws.Cells.Item(14, 8).Value = 281.031075
$ws.Cells.Item(14, 9).Value = 0.01594250007799006
$ws.Cells.Item(14, 10).Value = 0.01594250007799006
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 417.3096593333333
$ws.Cells.Item(14, 14).Value = 1251.928978
$ws.Cells.Item(14, 15).Value = 0.8277575129472603
$ws.Cells.Item(14, 16).Value = 0.8277575129472603
$ws.Cells.Item(14, 17).Value = 39092.32739011015
$ws.Cells.Item(14, 18).Value = 351830.9465109914
$ws.Cells.Item(14, 19).Value = 0.01319652421471856
$ws.Cells.Item(14, 20).Value = 0.01319652421471856
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 93.67702500000001
$ws.Cells.Item(15, 8).Value = 281.031075
$ws.Cells.Item(15, 9).Value = 0.01594250007799006
$ws.Cells.Item(15, 10).Value = 0.01594250007799006
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.194815333333333
$ws.Cells.Item(15, 14).Value = 3.584446
$ws.Cells.Item(15, 15).Value = 0.002369984366839822
$ws.Cells.Item(15, 16).Value = 0.002369984366839822
$ws.Cells.Item(15, 17).Value = 111.92674585105
$ws.Cells.Item(15, 18).Value = 1007.34071265945
$ws.Cells.Item(15, 19).Value = 0.00003778347595317909
$ws.Cells.Item(15, 20).Value = 0.00003778347595317908
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 93.67702500000001
$ws.Cells.Item(16, 8).Value = 281.031075
$ws.Cells.Item(16, 9).Value = 0.01594250007799006
$ws.Cells.Item(16, 10).Value = 0.01594250007799006
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 58.96764933333333
$ws.Cells.Item(16, 14).Value = 176.902948
$ws.Cells.Item(16, 15).Value = 0.1169656960121252
$ws.Cells.Item(16, 16).Value = 0.1169656960121252
$ws.Cells.Item(16, 17).Value = 5523.913960789901
$ws.Cells.Item(16, 18).Value = 49715.22564710911
$ws.Cells.Item(16, 19).Value = 0.001864725617795467
$ws.Cells.Item(16, 20).Value = 0.001864725617795467
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 93.67702500000001
$ws.Cells.Item(17, 8).Value = 281.031075
$ws.Cells.Item(17, 9).Value = 0.01594250007799006
$ws.Cells.Item(17, 10).Value = 0.01594250007799006
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 26.67269233333333
$ws.Cells.Item(17, 14).Value = 80.018077
$ws.Cells.Item(17, 15).Value = 0.05290680667377473
$ws.Cells.Item(17, 16).Value = 0.05290680667377473
$ws.Cells.Item(17, 17).Value = 2498.618466526975
$ws.Cells.Item(17, 18).Value = 22487.56619874278
$ws.Cells.Item(17, 19).Value = 0.0008434667695228587
$ws.Cells.Item(17, 20).Value = 0.0008434667695228587
